$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value  = 12.85
$ws.Range("E7").Value  = 13.044
$ws.Range("B8").Value  = 6.058999999999999
$ws.Range("A12").Value = -21.401
$ws.Range("B12").Value = 6.694999999999999
$ws.Range("B14").Value = 6.532000000000001
$ws.Range("E19").Value = 12.516
$ws.Range("E21").Value = 12.923
$ws.Range("B22").Value = 6.523999999999999
$ws.Range("E24").Value = 12.895
